# Update gh-pages output: add new event "广西·首届明日方舟only展 - 花庭圣梦" (2024-07-14)
# to the 展览 (Sheet1) and 全部类型 (Sheet4) sheets, and refresh the
# "想去人数" (want-to-go) counters that ticked up for several existing rows.

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $addr, $text) {
    # Force the cell to stay a literal text value (Excel otherwise
    # auto-parses strings like "2024-07-14" into date serials).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

function Add-NewRow($ws, $rowNum, $serial, $date, $title, $place, $timeRange, $wantCount, $price, $link, $cover) {
    $ws.Rows.Item($rowNum).Insert()

    $aAddr = "A" + $rowNum
    $ws.Range($aAddr).Value = $serial
    # Match the bordered/bold/centered style used by every other A-column cell.
    $ws.Range($aAddr).Borders.LineStyle = 1

    Set-TextValue $ws ("B" + $rowNum) $date

    $ws.Range("C" + $rowNum).Value = $title
    $ws.Range("D" + $rowNum).Value = $place
    $ws.Range("E" + $rowNum).Value = $timeRange
    $ws.Range("F" + $rowNum).Value = $wantCount
    $ws.Range("G" + $rowNum).Value = $price
    $ws.Range("H" + $rowNum).Value = $link
    $ws.Range("I" + $rowNum).Value = $cover
}

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value = 272
$ws1.Range("F4").Value = 27
$ws1.Range("F5").Value = 3307
$ws1.Range("F6").Value = 2119
$ws1.Range("F7").Value = 404
$ws1.Range("F8").Value = 153
$ws1.Range("F9").Value = 32

Add-NewRow $ws1 10 9 "2024-07-14" "广西·首届明日方舟only展 - 花庭圣梦" `
    "明秀东路157号 利泰国际大酒店" "2024.07.14 09:00-07.14 18:00" 9 69 `
    "https://show.bilibili.com/platform/detail.html?id=85852" `
    "//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg"

# Rows that slid down one slot keep their own data; only re-stamp the
# positional serial number in column A and the refreshed want-to-go counts.
$ws1.Range("A11").Value = 10
$ws1.Range("A12").Value = 11
$ws1.Range("A13").Value = 12
$ws1.Range("A14").Value = 13

$ws1.Range("F11").Value = 1213
$ws1.Range("F12").Value = 221
$ws1.Range("F13").Value = 1259
$ws1.Range("F14").Value = 102

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value = 272
$ws4.Range("F4").Value = 27
$ws4.Range("F5").Value = 3307
$ws4.Range("F6").Value = 2119
$ws4.Range("F7").Value = 404
$ws4.Range("F9").Value = 153
$ws4.Range("F10").Value = 32

Add-NewRow $ws4 11 10 "2024-07-14" "广西·首届明日方舟only展 - 花庭圣梦" `
    "明秀东路157号 利泰国际大酒店" "2024.07.14 09:00-07.14 18:00" 9 69 `
    "https://show.bilibili.com/platform/detail.html?id=85852" `
    "//i2.hdslb.com/bfs/openplatform/202405/xsMTmueN1715920435584.jpeg"

$ws4.Range("A12").Value = 11
$ws4.Range("A13").Value = 12
$ws4.Range("A14").Value = 13
$ws4.Range("A15").Value = 14
$ws4.Range("A16").Value = 15
$ws4.Range("A17").Value = 16

$ws4.Range("F14").Value = 1213
$ws4.Range("F15").Value = 221
$ws4.Range("F16").Value = 1259
$ws4.Range("F17").Value = 102
